$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Notes section (A13:J14, merged) with Saad's estimate note.
# Do this before changing the title so that new shared strings are appended
# in the same order as the target workbook (Notes text first, then title).
$ws.Range("A13").Value = "Saad thinks that he will be able to complete T15 in 3 story points instead of 5 as originally estimated"

# Center the notes text vertically as well as horizontally, and keep wrap text on.
$ws.Range("A13:J14").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A13:J14").VerticalAlignment = -4108     # xlCenter
$ws.Range("A13:J14").WrapText = $true

# Update the title banner from the old sprint 3 dates to the new sprint 4 dates.
$ws.Range("A1").Value = "Sprint 4 Plan (Nov 6 - Nov 10)"

# Mark the U4 and U9 user stories as Completed (column J holds the "X" marker).
$ws.Range("J4").Value = "X"
$ws.Range("J9").Value = "X"

# Move the active selection from J7 to the header row A2:J2.
$ws.Range("A2:J2").Select()
